# Apply "weight decay modified" edit to Observation.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25: fill in the Train/Test Final Accuracy values that were missing
$ws.Range("D25").Value = 78.16
$ws.Range("E25").Value = 68.51

# Row 26: new entry for an additional weight-decay run (mirrors row 25's
# "optimizer - weight decay" parameter, same shared-string / formatting)
$ws.Range("A26").Value = 17
$ws.Range("B26").Value = "optimizer - weight decay"
$ws.Range("C26").Value = 0.001

# Move the active selection to C26 to match the saved cursor position
$ws.Range("C26").Select()
